# Appends two new data rows (74 and 75) to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74: new contract entry (Usina Mixx Produções LTDA)
$ws.Range("K74").Value = "082/2024"
$ws.Range("L74").Value = "IN020-2024SEMCULTE"
$ws.Range("M74").Value = "27/06/2024"
$ws.Range("N74").Value = "Não informado"
$ws.Range("O74").Value = "Usina Mixx Produções LTDA"
$ws.Range("P74").Value = "49.607.556/0001-30"
$ws.Range("Q74").Value = "INEXIGIBILIDADE DE`nLICITAÇÃO"
$ws.Range("R74").Value = "Contratação de pessoa jurídica para realização de`nshow musical do(a) Artista MARLUS VIANA para apresentação nos Festejos de São Pedro 2024,`nno Município de Nilo Peçanha BA, conforme grade especificado pela Secretaria Municipal de`nTurismo, Cultura, Esporte e Lazer."
$ws.Range("S74").Value = "MUNICÍPIO DE NILO PEÇANHA"
$ws.Range("T74").Value = "130.000,00"

# Row 75: new contract entry (LABORCOM COMÉRCIO DE MATERIAIS DE CONSTRUÇÃO LTDA.)
$ws.Range("K75").Value = "154/2023"
$ws.Range("L75").Value = "068/2023"
$ws.Range("M75").Value = "19/12/2023"
$ws.Range("N75").Value = "31/12/2023"
$ws.Range("O75").Value = "LABORCOM COMÉRCIO DE MATERIAIS DE CONSTRUÇÃO LTDA."
$ws.Range("P75").Value = "34.101.659/0001-56"
$ws.Range("Q75").Value = "DISPENSA DE LICITAÇÃO"
$ws.Range("R75").Value = "Contratação de pessoa jurídica para fornecimento de material elétrico para ILUMINAÇÃO PÚBLICA em atendimento às necessidades da Secretaria de Infraestrutura e Urbanismo do Município de Nilo Peçanha - BA."
$ws.Range("S75").Value = "MUNICÍPIO DE NILO PEÇANHA"
$ws.Range("T75").Value = "54.720,00"
